# Nesto_TestCases.xlsx - "Tweleth commit with Pdf Download Test Implimented"
#
# Replace the old "Verify Report Download Button" test case (TC_DB_05's single
# verification step) with a new "Test Sales Report Download" test case that
# clicks the download button and verifies a "pdf" result, pushing the
# existing TC_DB_06 (Verify Logout Functionality) test case down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "Dashboard_Tests"

# Insert a new blank row at row 15. This shifts the old row 15
# (TC_DB_06 / Verify Logout Functionality / step 1) down to row 16, and the
# old row 16 (blank / blank / step 2 "Verify URL contains login") down to
# row 17 -- exactly matching the target layout without touching them further.
$ws.Rows.Item(15).Insert()

# --- Row 14: turn it into the new "Test Sales Report Download" header row ---
$ws.Range("B14").Value = "Test Sales Report Download"
$ws.Range("C14").Value = "1.Click on ""Download Button"" at ""//a[contains(@class, 'btn-report')]"""

# B14 previously had the default (no) style; give it the same boxed style
# used throughout column B/A in this table (style used by A14).
$ws.Range("A14").Copy()
$ws.Range("B14").PasteSpecial(-4122)   # xlPasteFormats

# C14 becomes a "header" step cell like C10/C12 (boxed + shaded fill).
$ws.Range("C10").Copy()
$ws.Range("C14").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 15 (new): second step of the new test case ---
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C15:J15").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C15").Value = "2.Verify ""pdf"""

# Match the saved selection/active cell shown in the target workbook.
$null = $ws.Range("C15").Select()
